$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D (lesion_volume), shifting FPM, WPM, WPM_log, FPM_log left by one column
$ws.Range("D1").EntireColumn.Delete()
